$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new column D (Categoría 3) values, header first, then the
# "Económicamente Activo" group before the "Inactivo" group so the
# shared-string table is rebuilt in the expected order.
$ws.Range("D1").Value = "Categoría 3"
$ws.Range("D5").Value = "Económicamente Activo"
$ws.Range("D6").Value = "Económicamente Activo"
$ws.Range("D7").Value = "Económicamente Activo"
$ws.Range("D8").Value = "Económicamente Activo"
$ws.Range("D9").Value = "Económicamente Activo"
$ws.Range("D10").Value = "Económicamente Activo"
$ws.Range("D11").Value = "Económicamente Activo"
$ws.Range("D12").Value = "Económicamente Activo"
$ws.Range("D13").Value = "Económicamente Activo"
$ws.Range("D14").Value = "Económicamente Activo"
$ws.Range("D2").Value = "Inactivo"
$ws.Range("D3").Value = "Inactivo"
$ws.Range("D4").Value = "Inactivo"
$ws.Range("D15").Value = "Inactivo"
$ws.Range("D16").Value = "Inactivo"
$ws.Range("D17").Value = "Inactivo"
$ws.Range("D18").Value = "Inactivo"

# Fix existing column C text for rows 5 and 6 (University education age range)
$ws.Range("C5").Value = "Educación Universitaria (16-25 años)"
$ws.Range("C6").Value = "Educación Universitaria (16-25 años)"

# Set column D width close to the original author's best-fit width
$ws.Columns.Item(4).ColumnWidth = 20
